$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Data" column (D), shifting
# Data/Note/USR one column to the right (D->E, E->F, F->G).
$ws.Columns("D:D").Insert()

# New header for the inserted column, matching the bold/wrap style used
# for the other header cells (same formatting as A1).
$ws.Range("D1").Value = "Label on`ngateway"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").WrapText = $true

# New "Label on gateway" values for each data row, with no special
# formatting (plain/default cell style).
$ws.Range("D4").Value = "V.1.010"
$ws.Range("D5").Value = "V.1.010"
$ws.Range("D6").Value = "V.1.020"

# Column D should keep the same width as column C.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth

# Row 6 FW VER value changes from "091" to the new "092" (kept as text).
$ws.Range("A6").Value = "'092"

$ws.Range("D9").Select()
